{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the text revisions described in the commit \"updated spec and UI\".\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// Each entry: literal text to find (must match exactly once) and its replacement.\nconst replacements = [\n  {\n    find: \"To display information about a module found in a .dll file and its associated .xml file in a TreeView \",\n    replace: \"To display information about all modules found in a directory using its dll and xml files in a TreeView \",\n  },\n  {\n    find: \"Under the module name, the public method names in that module will be displayed. Each method will be able to be selected which will then display all information about that module and the selected method.\",\n    replace: \"All public methods of a module will be displayed under the module name in the TreeView. The user can select each method which will then display all information about that module and the selected method in a different field.\",\n  },\n  {\n    find: \" own enable/disable check box to indicate whether the module is being used or not. \",\n    replace: \" own enable/disable check box to indicate whether the module will be used or not. \",\n  },\n  {\n    find: \"It will have a \\u201cload module location\\u201d button to allow the user to select the location of the dll\\u2019s.\",\n    replace: \"It will have a \\u201cload modules\\u201d button to allow the user to select the location of the dll\\u2019s.\",\n  },\n  {\n    find: \"The state of each enable/disable will be able to be saved to a config file. On startup of the application, this config file, if it exists, will be loaded and the enable/disable status of each module located in the previously specified file path will be displayed. \",\n    replace: \"There will be a \\u201csave configuration\\u201d button which will save the state of each modules enable/disable check box. On application startup, if the config file exists, the previously saved settings will be loaded. The last selected directory will also be saved so all modules in that directory will be loaded and displayed. \",\n  },\n  {\n    find: \"1.)  A load button will allow the user to select a directory where dll files are contained. \",\n    replace: \"1.)  A \\u201cload modules\\u201d button will allow the user to select a directory where dll files are located. The program will check that the directory contains at least one dll file.\",\n  },\n  {\n    find: \"3.)  The user can select any of the methods from the previously stated TreeView area. The information about the selected method will be displayed in a list box next to the TreeView area. Alternatively, this information can be displayed in a \\u201chover over\\u201d help text way.\",\n    replace: \"3.)  The user can select any of the methods from the previously stated TreeView area. The information about the selected method will be displayed in a separate field next to the TreeView area.\",\n  },\n  {\n    find: \", the application will load with this saved data as it exists. \",\n    replace: \", the application will load with this saved data. \",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the text revisions described in the commit \"updated spec and UI\".\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $found = $find.Execute(\n        [ref]$findText,    # FindText\n        [ref]$false,       # MatchCase\n        [ref]$false,       # MatchWholeWord\n        [ref]$false,       # MatchWildcards\n        [ref]$false,       # MatchSoundsLike\n        [ref]$false,       # MatchAllWordForms\n        [ref]$true,        # Forward\n        [ref]1,            # Wrap (wdFindContinue)\n        [ref]$false,       # Format\n        [ref]$replaceText, # ReplaceWith\n        [ref]2             # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-DocText `\n    \"To display information about a module found in a .dll file and its associated .xml file in a TreeView \" `\n    \"To display information about all modules found in a directory using its dll and xml files in a TreeView \"\n\nReplace-DocText `\n    \"Under the module name, the public method names in that module will be displayed. Each method will be able to be selected which will then display all information about that module and the selected method.\" `\n    \"All public methods of a module will be displayed under the module name in the TreeView. The user can select each method which will then display all information about that module and the selected method in a different field.\"\n\nReplace-DocText `\n    \" own enable/disable check box to indicate whether the module is being used or not. \" `\n    \" own enable/disable check box to indicate whether the module will be used or not. \"\n\nReplace-DocText `\n    \"It will have a \u201cload module location\u201d button to allow the user to select the location of the dll\u2019s.\" `\n    \"It will have a \u201cload modules\u201d button to allow the user to select the location of the dll\u2019s.\"\n\nReplace-DocText `\n    \"The state of each enable/disable will be able to be saved to a config file. On startup of the application, this config file, if it exists, will be loaded and the enable/disable status of each module located in the previously specified file path will be displayed. \" `\n    \"There will be a \u201csave configuration\u201d button which will save the state of each modules enable/disable check box. On application startup, if the config file exists, the previously saved settings will be loaded. The last selected directory will also be saved so all modules in that directory will be loaded and displayed. \"\n\nReplace-DocText `\n    \"1.)  A load button will allow the user to select a directory where dll files are contained. \" `\n    \"1.)  A \u201cload modules\u201d button will allow the user to select a directory where dll files are located. The program will check that the directory contains at least one dll file.\"\n\nReplace-DocText `\n    \"3.)  The user can select any of the methods from the previously stated TreeView area. The information about the selected method will be displayed in a list box next to the TreeView area. Alternatively, this information can be displayed in a \u201chover over\u201d help text way.\" `\n    \"3.)  The user can select any of the methods from the previously stated TreeView area. The information about the selected method will be displayed in a separate field next to the TreeView area.\"\n\nReplace-DocText `\n    \", the application will load with this saved data as it exists. \" `\n    \", the application will load with this saved data. \"\n"}
